$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values of row 2 and row 3 for columns D, M, N, O, P, S.
$cols = @("D", "M", "N", "O", "P", "S")

foreach ($col in $cols) {
    $cellRow2 = $ws.Range($col + "2")
    $cellRow3 = $ws.Range($col + "3")

    $v2 = $cellRow2.Value()
    $v3 = $cellRow3.Value()

    $cellRow2.Value = $v3
    $cellRow3.Value = $v2
}
